# "add two more lines"
#
# Before:
#   P1: "Hey hey"                       (empty "_GoBack" bookmark lives in P2)
#   P2: ""                              (bookmarkStart/bookmarkEnd "_GoBack")
#
# After:
#   P1: "August slipped away "          (bookmarkStart/bookmarkEnd "_GoBack" now leads P1)
#   P2: "Into moment in time"
#   P3: "Cause it was never mine"
#   P4: ""                              (the old, now-vacated, trailing paragraph)

$d = $word.ActiveDocument

# 1. Reword the first paragraph: "Hey hey" -> "August slipped away " (note the
#    trailing space kept from the source diff).
$d.Paragraphs(1).Range.Find.Execute("Hey hey", $false, $false, $false, $false, $false, $true, 1, $false, "August slipped away ", 2) | Out-Null

# 2. Insert two new paragraphs right after paragraph 1 (and therefore before the
#    original trailing paragraph that carries the "_GoBack" bookmark), then give
#    them the new lyric lines.
$p1 = $d.Paragraphs(1).Range
$p1.InsertParagraphAfter()
$p1.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "Into moment in time"
$d.Paragraphs(3).Range.Text = "Cause it was never mine"

# 3. Relocate the hidden "_GoBack" bookmark so it collapses at the very start of
#    paragraph 1 instead of sitting in the (now fourth) trailing empty paragraph.
#    A scratch leading paragraph is inserted/removed around the Add() call so the
#    bookmark can land exactly at document position 0.
$d.Paragraphs(1).Range.InsertParagraphBefore()
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range(1, 1))
$d.Paragraphs(1).Range.Delete()
